$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete years 2000年-2009年 (original rows 2-11),
# shifting 2010年-2020年 up to rows 2-12
$ws.Range("A2:A11").EntireRow.Delete()

# Append the new 2021年 row of data as row 13
# Copy formatting from the row above (A12) so the new row matches
# the existing bold/centered/bordered style used throughout column A
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 196.80677
$ws.Range("C13").Value = 655
$ws.Range("D13").Value = 443.6291026896
$ws.Range("E13").Value = 1422.54925442
